$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.093714853270462584
    2  = -0.0059999999140387672
    3  = -0.0039999999311568502
    4  = -0.0079999998724549215
    5  = -0.0029999999402230415
    6  = -0.0019999999482021025
    7  = -0.0099999998350721597
    8  = -0.0099999998352560127
    9  = -0.0019999999509172639
    10 = -0.0019999999543767188
    11 = -0.0029999999403882427
    12 = -0.0034999999352622879
    13 = 0.0042319541705655439
    14 = -0.0079999998831894459
    15 = 0.03411309792077688
    16 = -0.0019999999703741445
    17 = -0.0019999999686026726
    18 = -0.0039999999400555097
    19 = -0.0039999999419708665
    20 = -0.003999999940768717
    21 = -0.003999999941082244
    22 = 0.0095951682939476868
    23 = -0.0049999999115737381
    24 = -0.01999999968994981
    25 = -0.019999999685413883
    26 = -0.0024999999408592544
    27 = -0.002499999940439146
    28 = -0.0019999999452551265
    29 = -0.0069999998737273472
    30 = -0.059999999126776338
    31 = 0.01647415548313802
    32 = 0.050146151347254531
    33 = -0.0039999999357274163
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
